$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (zkratka) swaps: KEIS/EIS, KSWI/SWI, OPRE/ITM/KOPRE rotation ---
$ws.Cells.Item(2, 2).Value = "EIS"
$ws.Cells.Item(3, 2).Value = "KEIS"
$ws.Cells.Item(4, 2).Value = "EIS"
$ws.Cells.Item(5, 2).Value = "KEIS"
$ws.Cells.Item(6, 2).Value = "EIS"
$ws.Cells.Item(7, 2).Value = "KEIS"
$ws.Cells.Item(8, 2).Value = "EIS"
$ws.Cells.Item(9, 2).Value = "KEIS"
$ws.Cells.Item(14, 2).Value = "EIS"
$ws.Cells.Item(15, 2).Value = "KEIS"
$ws.Cells.Item(16, 2).Value = "SWI"
$ws.Cells.Item(17, 2).Value = "KSWI"

# --- Row 18/19: nazev (A) and zkratka (B) swap between "Odborna prezentace" and "Introduction to MATLAB" ---
$ws.Cells.Item(18, 1).Value = "Introduction to MATLAB"
$ws.Cells.Item(18, 2).Value = "ITM"
$ws.Cells.Item(19, 1).Value = "Odborná prezentace"
$ws.Cells.Item(19, 2).Value = "KOPRE"
$ws.Cells.Item(20, 2).Value = "OPRE"

# --- Column C (seminariciUcitIdno): convert text values to real numbers ---
$ws.Cells.Item(2, 3).Value = 14
$ws.Cells.Item(3, 3).Value = 14
$ws.Cells.Item(4, 3).Value = 1609
$ws.Cells.Item(5, 3).Value = 1609
$ws.Cells.Item(6, 3).Value = 3457
$ws.Cells.Item(7, 3).Value = 3457
$ws.Cells.Item(8, 3).Value = 3606
$ws.Cells.Item(9, 3).Value = 3606
$ws.Cells.Item(10, 3).Value = 4190
$ws.Cells.Item(11, 3).Value = 4303
$ws.Cells.Item(12, 3).Value = 4746
$ws.Cells.Item(13, 3).Value = 4746
$ws.Cells.Item(14, 3).Value = 4991
$ws.Cells.Item(15, 3).Value = 4991
$ws.Cells.Item(16, 3).Value = 8093
$ws.Cells.Item(17, 3).Value = 8093
$ws.Cells.Item(18, 3).Value = 8514
$ws.Cells.Item(19, 3).Value = 8514
$ws.Cells.Item(20, 3).Value = 8514

# Apply the new number format (thousands separator, red negatives) to the numeric id column
$ws.Range("C2:C20").NumberFormat = "#,##0;[Red]-#,##0"
